# database.xlsx -- "D-funcs: fix handling of blanks. Also update expected
# values.  We still fail." -- re-pointed at the Operators sheet (the D-func
# regression tests live there) and refreshed the hard-coded "expected
# result" fixture values to match the corrected behaviour.

$wb = $excel.ActiveWorkbook

$overview  = $wb.Worksheets.Item("Overview")
$types     = $wb.Worksheets.Item("Types")
$operators = $wb.Worksheets.Item("Operators")

# --- Shared-string edits -------------------------------------------------
# The "Regexp" label that used to be shared by the three wildcard tests now
# only applies to the ones that really are wildcard/regexp tests; give it a
# more precise name there...
$operators.Range("F48").Value = "Regexps only match strings"
$operators.Range("F52").Value = "Regexps only match strings"
# ...while the plain-string comparison test (row 80) keeps the original,
# simpler "Regexp" label.
$operators.Range("F80").Value = "Regexp"

# --- Operators!C44 / C48 / C52 -------------------------------------------
# These DSUM() calls now return different results once the blanks handling
# is fixed; record the corrected values the same way the fixture's other
# "actual result" cells are recorded (as plain numbers).
$operators.Range("C44").Value = 32751
$operators.Range("C48").Value = 0
$operators.Range("C52").Value = 16384

# --- Operators!C56 --------------------------------------------------------
# This DSUM() now matches its expected value (D56 = 2048), so the
# pass/fail check in E56 flips from FALSE to TRUE.
$operators.Range("C56").Value = 2048

# --- Operators!B1 / B2 (Pass / Fail totals) -------------------------------
# These totals roll up the E-column pass/fail flags for the whole sheet;
# refresh them to the corrected counts (15/1 -> 13/3).
$operators.Range("B1").Value = 13
$operators.Range("B2").Value = 3

# Overview!B3 / C3 reference Operators!B1 / Operators!B2 and recalculate on
# their own once the above are updated.

# --- Active sheet / selection --------------------------------------------
# The workbook used to open on Overview (tab 0); it now opens on Operators
# (tab 2), with the selection sitting on A1 there instead of B3.
$operators.Activate()
$operators.Range("A1").Select()
